# Daily attendance processing - reorder "Recorded By" (column G) entries
# The "System" / "system" token (and similarly admin@admin.com) is moved
# to a different position within the comma-separated list of recorders.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Exact value-level replacements observed in the target diff.
$replacements = @{
    "System, backup@backdoor.com"         = "backup@backdoor.com, System"
    "System, dnasr281@gmail.com"          = "dnasr281@gmail.com, System"
    "admin@admin.com, dnasr281@gmail.com" = "dnasr281@gmail.com, admin@admin.com"
    "system, System, backup@backdoor.com" = "system, backup@backdoor.com, System"
}

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = 7 ("Recorded By")
    $val = $cell.Value2
    if ($null -ne $val -and $replacements.ContainsKey($val)) {
        $cell.Value2 = $replacements[$val]
    }
}
